$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data row 2 values to match the new export snapshot
$ws.Range("A2").Value = "gWePH458"
$ws.Range("B2").Value = 23111506
$ws.Range("C2").Value = "getwnor29"
$ws.Range("D2").Value = "Xve%2$7N"
$ws.Range("F2").Value = "iBaNCJva"
$ws.Range("G2").Value = "veUY"
